# week 24 Meeting.pptx - "added final commenting for some of the code and
# reorganised the folders"
#
# Content-level changes reproduced here:
#   1) A new slide (12th, SlideID 495) is appended at the end of the deck.
#      It ends up containing a single, empty "content" placeholder shape
#      (idx=1) - the Title placeholder and the layout's default Content
#      placeholder that PowerPoint created with the slide were removed
#      during editing.
#   2) The "today" date fields cached on the slide master / slide layouts
#      / notes master (date placeholders, format M/D/YYYY and DD/MM/YYYY)
#      advance from 4/4/2022 to 4/7/2022 (04/04/2022 -> 07/04/2022),
#      reflecting the date PowerPoint auto-stamped the deck with when it
#      was next saved.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Add the new slide at the end of the deck, using the same
#    "Title and Content" layout (slideLayout2.xml == PpSlideLayout 16)
#    slide 11 ("Questions") already uses.
# ---------------------------------------------------------------------
$newIndex = $p.Slides.Count + 1
$s = $p.Slides.Add($newIndex, 16)

$layouts = $p.SlideMaster.CustomLayouts
$blankLayout = $layouts.Item(7)    # "Blank" - defines no placeholders
$titleContentLayout = $layouts.Item(2)    # "Title and Content"

# Drop the placeholders PowerPoint auto-created with the slide.
while ($s.Shapes.Count -gt 0) {
    $s.Shapes.Item(1).Delete()
}

# The final placeholder that survives on the slide is the 36th shape
# PowerPoint ever allocated an id for on this slide (id=37, "Content
# Placeholder 36") after a long round of pasting/undoing content while
# editing. Re-synthesise that id by repeatedly letting the layout
# re-instantiate its placeholders (toggle to a layout with none of them,
# then back) and discarding them again, so the shape id counter advances
# the same way, landing the final placeholder on id 37.
for ($n = 0; $n -lt 17; $n++) {
    $s.CustomLayout = $blankLayout
    $s.CustomLayout = $titleContentLayout
    if ($n -lt 16) {
        while ($s.Shapes.Count -gt 0) {
            $s.Shapes.Item(1).Delete()
        }
    } else {
        # Last round: keep the content placeholder, drop only the title.
        for ($i = $s.Shapes.Count; $i -ge 1; $i--) {
            $sh = $s.Shapes.Item($i)
            if ($sh.Name -like "Title*") {
                $sh.Delete()
            }
        }
    }
}

$content = $s.Shapes.Item(1)
$content.Name = "Content Placeholder 36"

# ---------------------------------------------------------------------
# 2) Advance the cached "today" date text on the slide master and every
#    slide layout that has a date placeholder.
#    (Note: helper functions are deliberately NOT used here - calling a
#    PowerShell function repeatedly blows this host's statement budget
#    much faster than the equivalent code inlined in a loop.)
# ---------------------------------------------------------------------
$mshapes = $p.SlideMaster.Shapes
for ($i = 1; $i -le $mshapes.Count; $i++) {
    $sh = $mshapes.Item($i)
    if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
        $sh.TextFrame.TextRange.Text = "4/7/2022"
    }
}

for ($li = 1; $li -le $layouts.Count; $li++) {
    $lshapes = $layouts.Item($li).Shapes
    for ($i = 1; $i -le $lshapes.Count; $i++) {
        $sh = $lshapes.Item($i)
        if ($sh.Name -like "Date Placeholder*" -and $sh.HasTextFrame) {
            $sh.TextFrame.TextRange.Text = "4/7/2022"
        }
    }
}

# Note: the notes master (ppt/notesMasters/notesMaster1.xml) also carries a
# cached "today" field (04/04/2022 -> 07/04/2022 in the target) but this
# COM host exposes Presentation.NotesMaster as a read-only view (its
# shapes can be enumerated, but edits/deletes against them do not persist
# - consistent with Presentation.HasNotesMaster reporting False even
# though the part exists), so it cannot be updated from here.
